# Apply the "Finished initial part with data" edit to DF.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Fill in previously-empty ("inline string") cells in existing
#    data rows with the numeric value 0.
# ---------------------------------------------------------------
$zeroFixups = @(
    "J5", "K5", "T5", "W5", "AK5", "AM5",
    "D7",
    "AC10", "AN10"
)
foreach ($cellRef in $zeroFixups) {
    $ws.Range($cellRef).Value = 0
}

# ---------------------------------------------------------------
# 2. Append a new "World" aggregate row (row 12) under the last
#    existing data row (row 11), carrying over the bold/bordered
#    style used by the other rows in column A.
# ---------------------------------------------------------------
$ws.Range("A11").Copy($ws.Range("A12"))
$ws.Range("A12").Value = 213
$ws.Range("B12").Value = "World"

$row12Values = @{
    "C12"  = 88.84863805865645
    "D12"  = 78.64735578125527
    "E12"  = 97.35769156507307
    "F12"  = 42.45569899944301
    "G12"  = 54.21683086088953
    "H12"  = 13.21604411832126
    "I12"  = 39.6440809675836
    "J12"  = 3.418853897824389
    "K12"  = 0.8038986377878641
    "L12"  = 18.65693389189641
    "M12"  = 7.54240423884828
    "N12"  = 28.40931040761491
    "O12"  = 23.07488912516291
    "P12"  = 48.51578710374548
    "Q12"  = 2.431584662931561
    "R12"  = 80891341462319.14
    "S12"  = 3.165402310413128
    "T12"  = 25.50045483934074
    "U12"  = 72.38300866126713
    "V12"  = 44.47336011088299
    "W12"  = 2.161438045828547
    "X12"  = 102.7759657782136
    "Y12"  = 0
    "Z12"  = 25.93487403210847
    "AA12" = 65.41928443397815
    "AB12" = 8.645841533913378
    "AC12" = 58.97696677271883
    "AD12" = 1.143092581314392
    "AE12" = 16.09383079924788
    "AF12" = 0.4750067445347668
    "AG12" = 3722940052
    "AH12" = 3785204133
    "AI12" = 7510990456
    "AJ12" = 0.8
    "AK12" = 10.8
    "AL12" = 45.17775825906573
    "AM12" = 65.03571648102344
    "AN12" = 132036620.8118786
    "AO12" = 54.8222417409343
    "AP12" = 1.98513545849886
}

foreach ($cellRef in $row12Values.Keys) {
    $ws.Range($cellRef).Value = $row12Values[$cellRef]
}

Write-Host "Applied DF.xlsx update: zero-filled blanks + appended World row 12"
